$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 102 (004339183 / JALISON / 95.69)
# so the new account (004231509 / THEOMAR / 98.64) is placed in correct
# descending-Saldo sort order between 99.18 and 95.69.
$ws.Rows.Item(102).Insert()

# Leading apostrophe forces Excel to store the account number as text so the
# leading zeros are preserved (otherwise "004231509" would be read as a number).
$ws.Cells.Item(102, 1).Value = "'004231509"
$ws.Cells.Item(102, 2).Value = "THEOMAR"
$ws.Cells.Item(102, 3).Value = 98.64

# The original THEOMAR row (004231509 / THEOMAR / -1.36) has shifted down one
# row because of the insert above; it is now on row 253. Replace it with the
# new account (004384258 / PAULA / -6.71).
$ws.Cells.Item(253, 1).Value = "'004384258"
$ws.Cells.Item(253, 2).Value = "PAULA"
$ws.Cells.Item(253, 3).Value = -6.71
